$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 3 ("Poblacion") -- pushes existing rows 3..35 down to 4..36
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "Poblacion"
$ws.Range("B3").Value = "poblacion"

# 2. Insert a new row at position 15 ("Jerarquia cat") -- pushes rows 15..36 down to 16..37
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "Jerarquia cat"
$ws.Range("B15").Value = "jerarquia_cat"

# 3. Append new rows 38..55 after "Prestador" (now row 37)
$newRows = @(
    @("Tipo de vivienda", "tipo_de_vivienda"),
    @("Tipo de instrumento", "tipo_de_instrumento"),
    @("Tipo de inspección", "tipo_de_inspeccion"),
    @("Tipo de centro", "tipo_de_centro"),
    @("Sector económico", "sector_economico"),
    @("Prevalencia de limitaciones", "prevalencia_de_limitaciones"),
    @("Motivo por el cual se fueron", "motivo_por_el_cual_se_fueron"),
    @("Identidad de género", "identidad_de_genero"),
    @("Área geográfica", "area_geografica"),
    @("Cuenca o embalse", "cuenca_o_embalse"),
    @("Río", "rio"),
    @("Área de capacitación", "area_de_capacitacion"),
    @("Amparo bse", "amparo_bse"),
    @("Condición migratoria", "condicion_migratoria"),
    @("Adecuación educativa", "adecuacion_educativa"),
    @("Tipo de cláusula", "tipo_de_clausula"),
    @("Causa de muerte", "causa_de_muerte"),
    @("Situación procesal", "situacion_procesal")
)

$row = 38
foreach ($pair in $newRows) {
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("B$row").Value = $pair[1]
    $row = $row + 1
}
